# Default language: Portuguese + smaller changes
#
# Swap the default display language from English to Portuguese across the
# "settings" and "choices" sheets (fixing the "tittle" typo along the way),
# and touch up the active-sheet / selection state to match the author's
# final editing session.

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")
$wsSurvey   = $wb.Worksheets.Item("survey")
$wsChoices  = $wb.Worksheets.Item("choices")

# -----------------------------------------------------------------------
# Header renames.
#
# The ".english" suffixed headers are brand new column names, so they are
# written first (and in this particular order) to line up with how the
# workbook's shared-string table ends up laid out. The plain (no-suffix)
# headers simply become aliases of the already-existing
# "display.title.text" / "display.locale.text" / "display.prompt.text"
# strings, now repurposed as the Portuguese/default columns.
# -----------------------------------------------------------------------

# settings sheet: title headers
$wsSettings.Range("C1").Value = "display.title.text.english"

# survey sheet: prompt headers
$wsSurvey.Range("G1").Value = "display.prompt.text.english"

# settings sheet: locale headers
$wsSettings.Range("E1").Value = "display.locale.text.english"

# settings sheet: the former "portuguese" settings row is renamed "english"
$wsSettings.Range("A7").Value = "english"

# settings sheet: plain headers now mean the (Portuguese) default column
$wsSettings.Range("D1").Value = "display.title.text"
$wsSettings.Range("F1").Value = "display.locale.text"

# survey sheet: plain header now means the (Portuguese) default column
$wsSurvey.Range("H1").Value = "display.prompt.text"

# choices sheet: same header rename as settings
$wsChoices.Range("C1").Value = "display.title.text.english"
$wsChoices.Range("D1").Value = "display.title.text"

# -----------------------------------------------------------------------
# Default-language row on the settings sheet: the "default" row now holds
# the Portuguese strings, and the (renamed) "english" row holds what used
# to be the default English strings.
# -----------------------------------------------------------------------
$wsSettings.Range("E6").Value = "Portuguese"
$wsSettings.Range("F6").Value = "Portugues"
$wsSettings.Range("E7").Value = "English"
$wsSettings.Range("F7").Value = "Inglês"

# -----------------------------------------------------------------------
# Selections / active sheet, matching the author's final cursor position.
#
# Selecting a range on a worksheet implicitly activates that worksheet, so
# these run in an order that leaves "settings" activated last (it is the
# tab that ends up selected in the saved workbook).
# -----------------------------------------------------------------------

# survey: selection moves to the edited header cells G1:H1.
$wsSurvey.Range("G1:H1").Select()

# choices: selection moves to the edited header cells C1:D1.
$wsChoices.Range("C1:D1").Select()

# settings becomes the active tab, with A6:F7 selected (the rows that were
# just edited).
$wsSettings.Activate()
$wsSettings.Range("A6:F7").Select()
